# "Add files via upload" — the author typed the real Centro Gestor unit
# names (previously abbreviated as "52010 UGA" / "52012 UPDN") into the
# Ingresos_Centros sheet, bumped the assigned-income figures, and added a
# new row for the fourth unit. The selections left behind in both sheets
# (Grupos_Centros B2:B5 — the source of the unit names — and
# Ingresos_Centros A6, the next empty row) are reproduced too.

$wb = $excel.ActiveWorkbook

$wsGrupos   = $wb.Worksheets.Item("Grupos_Centros")
$wsIngresos = $wb.Worksheets.Item("Ingresos_Centros")

# Update the existing rows with the full unit names and new amounts.
$wsIngresos.Range("B2").Value = 1000000

$wsIngresos.Range("A3").Value = "52010 UD GESTIÓN ADMINISTRATIVA"
$wsIngresos.Range("B3").Value = 800000

$wsIngresos.Range("A4").Value = "52011 UD GENERACIÓN ENERGIA"
$wsIngresos.Range("B4").Value = 2000000

# New fourth row.
$wsIngresos.Range("A5").Value = "52012 UD PROSPECTIVA Y DESARROLLO NEGOCIOS"
$wsIngresos.Range("B5").Value = 550000

# Column A now holds much longer labels -- widen it to fit (~45.4 chars).
$wsIngresos.Columns.Item(1).ColumnWidth = 44.67

# Leave behind the same selections the author had: the source range they
# copied the unit names from on Grupos_Centros, ...
[void]$wsGrupos.Range("B2:B5").Select()

# ... and the next blank row on Ingresos_Centros, which stays the active
# (tab-selected) sheet.
[void]$wsIngresos.Activate()
[void]$wsIngresos.Range("A6").Select()
